# Hoàn thiện Ngoại Trú
$wb = $excel.ActiveWorkbook

# --- "Data" sheet (sheet1) ---
$data = $wb.Worksheets.Item("Data")

$data.Range("A2").Value = 3005
$data.Range("E2").Value = 46200608005
$data.Range("W2").Value = "5/49 Ntl"
$data.Range("X2").Value = "DN4127460130005"
$data.Range("BF2").Value = "||1|Normal|CorrectRoute|09/05/2024 09:13|3266971|Quách Bảo Hưng 82|24|Male|01/01/2000 00:00|5/49 Ntl|765|26926|01|VN|134||DN4127389127512|2|80|None|||||||||3839|Open|||||||||||149|09/05/2024 09:13||||||3839||||New|4803|80|||2|1083660|||"

# Update the active selection to X2, matching the saved view state
$data.Activate()
$data.Range("X2").Select()

# --- "Check" sheet (sheet2) ---
$check = $wb.Worksheets.Item("Check")

$check.Range("A2").Value = 3005
$check.Range("C2").Value = "DN4127460130005"

Write-Host "Applied Hoan thien Ngoai Tru edits"
